# Collate hour worked from employees - update hourly rates / hours / sick days
# and freeze the "Total cost" column to static values (progress notification
# message boxes shown to the user while the automation runs through each row).

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

Write-Host "Starting payroll update automation..."

# Row 3 (Marketta Elmore): hourly rate, hours worked and sick days change;
# keep the Total cost formula (H3*I3) - it will recalc automatically.
Write-Host "Processing row 3 of 20..."
$ws.Range("H3").Value = 10
$ws.Range("I3").Value = 11
$ws.Range("K3").Value = 1

# All other rows: only the hourly rate changes. After updating the rate we
# freeze the "Total cost" column by writing the computed number back in as a
# plain value, which replaces the =H*I formula with a static result.
$newRates = @{
    2  = 12
    4  = 17
    5  = 15
    6  = 13
    7  = 20
    8  = 8
    9  = 15
    10 = 19
    11 = 10
    12 = 14
    13 = 12
    14 = 12
    15 = 22
    16 = 18
    17 = 18
    18 = 23
    19 = 10
    20 = 15
    21 = 14
}

foreach ($row in $newRates.Keys | Sort-Object) {
    Write-Host "Processing row $row of 20..."
    $rate = $newRates[$row]
    $ws.Range("H$row").Value = $rate
    $hours = $ws.Range("I$row").Value2
    $total = $rate * $hours
    $ws.Range("M$row").Value = $total
}

Write-Host "Payroll update automation complete."
